# Appends the Warriors' final four regular-season games (2025-04-11 vs POR
# x2, 2025-04-13 vs LAC x2) to the box-score team-stats sheet, continuing
# directly after the existing last row (161).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 162 (game index 160) ---
$ws.Range("A161").Copy($ws.Range("A162"))
$ws.Range("A162").Value = 160
$ws.Range("B162").Value = "GSW"
$ws.Range("C162").Value = "POR"
$ws.Range("D162").Value = "away"
$ws.Range("E162").NumberFormat = "@"
$ws.Range("E162").Value = "2025-04-11"
$ws.Range("E162").ClearFormats()
$ws.Range("F162").Value = "240:00"
$ws.Range("G162").Value = 33
$ws.Range("H162").Value = 83
$ws.Range("I162").Value = 0.398
$ws.Range("J162").Value = 14
$ws.Range("K162").Value = 44
$ws.Range("L162").Value = 0.318
$ws.Range("M162").Value = 23
$ws.Range("N162").Value = 28
$ws.Range("O162").Value = 0.821
$ws.Range("P162").Value = 7
$ws.Range("Q162").Value = 38
$ws.Range("R162").Value = 45
$ws.Range("S162").Value = 26
$ws.Range("T162").Value = 15
$ws.Range("U162").Value = 7
$ws.Range("V162").Value = 14
$ws.Range("W162").Value = 19
$ws.Range("X162").Value = 103
$ws.Range("Y162").Value = 17
$ws.Range("Z162").Value = 26
$ws.Range("AA162").Value = 24
$ws.Range("AB162").Value = 27
$ws.Range("AC162").Value = 26
$ws.Range("AD162").Value = "W"

# --- Row 163 (game index 161) ---
$ws.Range("A162").Copy($ws.Range("A163"))
$ws.Range("A163").Value = 161
$ws.Range("B163").Value = "POR"
$ws.Range("C163").Value = "GSW"
$ws.Range("D163").Value = "home"
$ws.Range("E163").NumberFormat = "@"
$ws.Range("E163").Value = "2025-04-11"
$ws.Range("E163").ClearFormats()
$ws.Range("F163").Value = "240:00"
$ws.Range("G163").Value = 32
$ws.Range("H163").Value = 85
$ws.Range("I163").Value = 0.376
$ws.Range("J163").Value = 6
$ws.Range("K163").Value = 30
$ws.Range("L163").Value = 0.2
$ws.Range("M163").Value = 16
$ws.Range("N163").Value = 22
$ws.Range("O163").Value = 0.727
$ws.Range("P163").Value = 13
$ws.Range("Q163").Value = 39
$ws.Range("R163").Value = 52
$ws.Range("S163").Value = 21
$ws.Range("T163").Value = 9
$ws.Range("U163").Value = 5
$ws.Range("V163").Value = 22
$ws.Range("W163").Value = 20
$ws.Range("X163").Value = 86
$ws.Range("Y163").Value = -17
$ws.Range("Z163").Value = 20
$ws.Range("AA163").Value = 17
$ws.Range("AB163").Value = 20
$ws.Range("AC163").Value = 29
$ws.Range("AD163").Value = "L"

# --- Row 164 (game index 162) ---
$ws.Range("A163").Copy($ws.Range("A164"))
$ws.Range("A164").Value = 162
$ws.Range("B164").Value = "LAC"
$ws.Range("C164").Value = "GSW"
$ws.Range("D164").Value = "away"
$ws.Range("E164").NumberFormat = "@"
$ws.Range("E164").Value = "2025-04-13"
$ws.Range("E164").ClearFormats()
$ws.Range("F164").Value = "265:00"
$ws.Range("G164").Value = 48
$ws.Range("H164").Value = 86
$ws.Range("I164").Value = 0.558
$ws.Range("J164").Value = 14
$ws.Range("K164").Value = 30
$ws.Range("L164").Value = 0.467
$ws.Range("M164").Value = 14
$ws.Range("N164").Value = 18
$ws.Range("O164").Value = 0.778
$ws.Range("P164").Value = 9
$ws.Range("Q164").Value = 33
$ws.Range("R164").Value = 42
$ws.Range("S164").Value = 28
$ws.Range("T164").Value = 11
$ws.Range("U164").Value = 3
$ws.Range("V164").Value = 16
$ws.Range("W164").Value = 21
$ws.Range("X164").Value = 124
$ws.Range("Y164").Value = 5
$ws.Range("Z164").Value = 25
$ws.Range("AA164").Value = 35
$ws.Range("AB164").Value = 23
$ws.Range("AC164").Value = 28
$ws.Range("AD164").Value = "W"

# --- Row 165 (game index 163) ---
$ws.Range("A164").Copy($ws.Range("A165"))
$ws.Range("A165").Value = 163
$ws.Range("B165").Value = "GSW"
$ws.Range("C165").Value = "LAC"
$ws.Range("D165").Value = "home"
$ws.Range("E165").NumberFormat = "@"
$ws.Range("E165").Value = "2025-04-13"
$ws.Range("E165").ClearFormats()
$ws.Range("F165").Value = "265:00"
$ws.Range("G165").Value = 43
$ws.Range("H165").Value = 79
$ws.Range("I165").Value = 0.544
$ws.Range("J165").Value = 15
$ws.Range("K165").Value = 33
$ws.Range("L165").Value = 0.455
$ws.Range("M165").Value = 18
$ws.Range("N165").Value = 23
$ws.Range("O165").Value = 0.783
$ws.Range("P165").Value = 3
$ws.Range("Q165").Value = 22
$ws.Range("R165").Value = 25
$ws.Range("S165").Value = 31
$ws.Range("T165").Value = 11
$ws.Range("U165").Value = 6
$ws.Range("V165").Value = 15
$ws.Range("W165").Value = 20
$ws.Range("X165").Value = 119
$ws.Range("Y165").Value = -5
$ws.Range("Z165").Value = 33
$ws.Range("AA165").Value = 25
$ws.Range("AB165").Value = 24
$ws.Range("AC165").Value = 29
$ws.Range("AD165").Value = "L"

